$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header strings for columns L through S (row 1)
$headers = @(
    "hzj-混合调节_20170516_152754_ASIC_EEG",
    "hzj-混合调节_20170518_134207_ASIC_EEG",
    "hzj-混合调节_20170519_135415_ASIC_EEG",
    "zyx-混合调节_20170516_111228_ASIC_EEG",
    "zyx-混合调节_20170517_110944_ASIC_EEG",
    "zyx-混合调节_20170518_112337_ASIC_EEG",
    "zyx-混合调节_20170519_124954_ASIC_EEG",
    "zyx-混合调节_20170522_111557_ASIC_EEG"
)

# New data values for row 2 and row 3, columns L through S
$row2 = @(0.5625, 0.59450171821305842, 0.62172284644194753, 0.59935897435897434, 0.70192307692307687, 0.69255663430420711, 0.6064516129032258, 0.65776699029126218)
$row3 = @(0.53479853479853479, 0.57741935483870965, 0.58608058608058611, 0.57857142857142851, 0.65363128491620115, 0.60824742268041243, 0.60738255033557054, 0.63203463203463206)

$startCol = 12  # column L
for ($i = 0; $i -lt $headers.Count; $i++) {
    $col = $startCol + $i
    $ws.Cells.Item(1, $col).Value = $headers[$i]
    $ws.Cells.Item(2, $col).Value = $row2[$i]
    $ws.Cells.Item(3, $col).Value = $row3[$i]
}

# Update selection to new extent A1:S3
$ws.Range("A1:S3").Select()

# Update workbook window height
$excel.ActiveWindow.Height = 11655
